$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A64").Value = "19-11-2025"
$ws.Range("B64").Value = "The price of gold in India today is ₹12,486 per gram for 24 karat gold, ₹11,445 per gram for 22 karat gold and ₹9,364 per gram for 18 karat gold (also called 999 gold)."
